# Publishing after approval to publish: converts the page's four
# heading paragraphs (Heading1/Heading2) into Markdown-style "#"/"##"
# prefixed lines using the body's FirstParagraph style, demotes the
# paragraph that used to follow each heading from FirstParagraph to
# BodyText, and drops the now-unused heading bookmarks.

$d = $word.ActiveDocument

function Find-ParaIndex {
    param([string]$Text)
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $Text) {
            return $i
        }
    }
    return -1
}

function Convert-Heading {
    param([string]$HeadingText, [string]$Prefix)

    $headingIndex = Find-ParaIndex $HeadingText
    $heading = $d.Paragraphs.Item($headingIndex)

    # Insert a new (still empty, plain-formatted) paragraph right after
    # the heading, give it the FirstParagraph style, and fill it with
    # the markdown-ified title.
    $heading.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($headingIndex + 1)
    $newPara.Style = "FirstParagraph"
    $newPara.Range.Text = $Prefix + $HeadingText

    # The paragraph that used to directly follow the heading (now two
    # slots down) drops from FirstParagraph down to BodyText.
    $demoted = $d.Paragraphs.Item($headingIndex + 2)
    $demoted.Style = "BodyText"

    # Remove the original Heading1/Heading2 paragraph entirely.
    $heading = $d.Paragraphs.Item($headingIndex)
    $heading.Range.Delete()
}

Convert-Heading "Avoiding too much security" "#"
Convert-Heading "Not all domain names or IP addresses in Government systems are sensitive items" "##"
Convert-Heading "It’s not only about domain names or IP addresses" "##"
Convert-Heading "Feedback" "##"
